# Add a new Space Wolves stratagem row ("Deed worthy saga") and fix the
# timing text of the "Only in death, duty ends" stratagem (row 9) to use
# the same HTML-entity style as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build row 10 by copying row 9's formatting (A9:G9, no SourceLogo) --
$ws.Range("A9:G9").Copy() | Out-Null
$ws.Range("A10:G10").PasteSpecial() | Out-Null

# --- Populate the new row's content (same write order the author used) --
$ws.Range("B10").Value = "Toutes phases"
$ws.Range("A10").Value = "Deed worthy saga"
$ws.Range("C10").Value = "Lorsqu'un personnage space wolf autre que le seigneur de guerre satisfait les pré requis d'une saga, celui ci b&eacute;n&eacute;ficie des effets de celle ci jusqu'&agrave; la fin de la saga."
$ws.Range("D10").Value = $ws.Range("D9").Value()
$ws.Range("E10").Value = $ws.Range("E9").Value()
$ws.Range("F10").Value = $ws.Range("F9").Value()
$ws.Range("G10").Value = $ws.Range("G9").Value()

# --- Fix row 9's Timing cell (B9): "détruite" -> "d&eacute;truite" -------
$ws.Range("B9").Value = "Phase de combat, une figurine est d&eacute;truite"

# --- Row height for the new row -------------------------------------------
$ws.Rows.Item(10).RowHeight = 84

# --- Selection, matching the authored file's final cursor position -------
$ws.Range("C15").Select() | Out-Null
